# rnaSample_hbrown_06.11.19.xlsx — "changing FALSE to False"
#
# Column H (roboticRNAPrep) held boolean FALSE values displayed via a
# custom "TRUE"/"FALSE" number format. Re-enter them as the literal text
# "False" (stored as a shared string) formatted with the built-in Text
# number format instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-format H2:H27 as Text (built-in format 49, "@") before typing the
# literal word so it is stored as a string rather than re-parsed as a
# boolean.
$ws.Range("H2:H27").NumberFormat = "@"

# Leading apostrophe forces literal text entry (otherwise "False" would be
# auto-recognised as the Boolean FALSE, same as before).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 8).Value = "'False"
}

# Match the author's final selection/view state.
$ws.Range("H3:H27").Select()
$excel.ActiveWindow.Height = 1400
